$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2178571428571429
$ws.Range("C2").Value = 0.4928571428571429
$ws.Range("J2").Value = 0.01071428571428571
$ws.Range("P2").Value = 0.1785714285714286
$ws.Range("S2").Value = 0.1
$ws.Range("B3").Value = 0.0131578947368421
$ws.Range("C3").Value = 0.03289473684210526
$ws.Range("J3").Value = 0.02631578947368421
$ws.Range("P3").Value = 0.743421052631579
$ws.Range("S3").Value = 0.1842105263157895
$ws.Range("P4").Value = 0.6486486486486487
$ws.Range("S4").Value = 0.3513513513513514
$ws.Range("B6").Value = 0.03097345132743363
$ws.Range("D6").Value = 0.01327433628318584
$ws.Range("F6").Value = 0.06194690265486726
$ws.Range("J6").Value = 0.247787610619469
$ws.Range("O6").Value = 0.004424778761061947
$ws.Range("Q6").Value = 0.1946902654867257
$ws.Range("R6").Value = 0.1017699115044248
$ws.Range("S6").Value = 0.3451327433628318
$ws.Range("B7").Value = 0.09554140127388536
$ws.Range("D7").Value = 0.02547770700636943
$ws.Range("E7").Value = 0.006369426751592357
$ws.Range("F7").Value = 0.04458598726114649
$ws.Range("J7").Value = 0.1273885350318471
$ws.Range("O7").Value = 0.006369426751592357
$ws.Range("Q7").Value = 0.178343949044586
$ws.Range("R7").Value = 0.08280254777070063
$ws.Range("S7").Value = 0.4331210191082803
$ws.Range("B8").Value = 0.08823529411764706
$ws.Range("D8").Value = 0.02036199095022624
$ws.Range("E8").Value = 0.002262443438914027
$ws.Range("F8").Value = 0.05656108597285068
$ws.Range("J8").Value = 0.08597285067873303
$ws.Range("O8").Value = 0.01131221719457014
$ws.Range("Q8").Value = 0.1855203619909502
$ws.Range("R8").Value = 0.08823529411764706
$ws.Range("S8").Value = 0.4615384615384616
$ws.Range("B9").Value = 0.1004016064257028
$ws.Range("D9").Value = 0.02008032128514056
$ws.Range("E9").Value = 0.004016064257028112
$ws.Range("F9").Value = 0.05220883534136546
$ws.Range("J9").Value = 0.07630522088353414
$ws.Range("Q9").Value = 0.2008032128514056
$ws.Range("R9").Value = 0.06827309236947791
$ws.Range("S9").Value = 0.4779116465863454
$ws.Range("B10").Value = 0.1012658227848101
$ws.Range("D10").Value = 0.0134493670886076
$ws.Range("E10").Value = 0.0007911392405063291
$ws.Range("F10").Value = 0.07041139240506329
$ws.Range("J10").Value = 0.09731012658227849
$ws.Range("O10").Value = 0.01977848101265823
$ws.Range("Q10").Value = 0.245253164556962
$ws.Range("R10").Value = 0.07832278481012658
$ws.Range("S10").Value = 0.3734177215189873
$ws.Range("F11").Value = 0.00390625
$ws.Range("J11").Value = 0.11328125
$ws.Range("K11").Value = 0.19921875
$ws.Range("L11").Value = 0.55078125
$ws.Range("S11").Value = 0.0078125
$ws.Range("G12").Value = 0.7533333333333333
$ws.Range("J12").Value = 0.1466666666666667
$ws.Range("K12").Value = 0.02666666666666667
$ws.Range("L12").Value = 0.05333333333333334
$ws.Range("S12").Value = 0.02
$ws.Range("G13").Value = 0.8148148148148148
$ws.Range("J13").Value = 0.1481481481481481
$ws.Range("S13").Value = 0.03703703703703703
$ws.Range("F15").Value = 0.02304147465437788
$ws.Range("H15").Value = 0.152073732718894
$ws.Range("I15").Value = 0.09216589861751152
$ws.Range("J15").Value = 0.3686635944700461
$ws.Range("K15").Value = 0.04147465437788019
$ws.Range("M15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.04608294930875576
$ws.Range("S15").Value = 0.271889400921659
$ws.Range("F16").Value = 0.03867403314917127
$ws.Range("H16").Value = 0.138121546961326
$ws.Range("I16").Value = 0.08287292817679558
$ws.Range("J16").Value = 0.4751381215469613
$ws.Range("K16").Value = 0.08839779005524862
$ws.Range("M16").Value = 0.02209944751381215
$ws.Range("N16").Value = 0.005524861878453038
$ws.Range("O16").Value = 0.05524861878453038
$ws.Range("S16").Value = 0.09392265193370165
$ws.Range("F17").Value = 0.02339181286549707
$ws.Range("H17").Value = 0.1695906432748538
$ws.Range("I17").Value = 0.1111111111111111
$ws.Range("J17").Value = 0.442495126705653
$ws.Range("K17").Value = 0.0682261208576998
$ws.Range("M17").Value = 0.009746588693957114
$ws.Range("O17").Value = 0.07797270955165692
$ws.Range("S17").Value = 0.09746588693957114
$ws.Range("F18").Value = 0.0425531914893617
$ws.Range("H18").Value = 0.1436170212765958
$ws.Range("I18").Value = 0.1170212765957447
$ws.Range("J18").Value = 0.4148936170212766
$ws.Range("K18").Value = 0.1117021276595745
$ws.Range("M18").Value = 0.02127659574468085
$ws.Range("N18").Value = 0.005319148936170213
$ws.Range("O18").Value = 0.05851063829787234
$ws.Range("S18").Value = 0.0851063829787234
$ws.Range("F19").Value = 0.01413982717989002
$ws.Range("H19").Value = 0.213668499607227
$ws.Range("I19").Value = 0.1060487038491752
$ws.Range("J19").Value = 0.3880597014925373
$ws.Range("K19").Value = 0.09112333071484682
$ws.Range("M19").Value = 0.01178318931657502
$ws.Range("N19").Value = 0.0007855459544383347
$ws.Range("O19").Value = 0.06991358994501179
$ws.Range("S19").Value = 0.1044776119402985
